$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New roster table (18 players), row 2..19, columns A=Oyuncu Adı, B=Pozisyon, C=Takım
$rows = @(
    @("Cade Cunningham", "PG,SG", "Detroit Pistons"),
    @("Derrick White", "PG,SG", "Boston Celtics"),
    @("Bilal Coulibaly", "SG,SF", "Washington Wizards"),
    @("Ben Simmons", "PG,C", "Brooklyn Nets"),
    @("Duncan Robinson", "SG,SF", "Miami Heat"),
    @("Cameron Johnson", "SF,PF", "Brooklyn Nets"),
    @("Herbert Jones", "SF,PF", "New Orleans Pelicans"),
    @("Anthony Davis", "PF,C", "Los Angeles Lakers"),
    @("Bam Adebayo", "C", "Miami Heat"),
    @("Julius Randle", "PF,C", "Minnesota Timberwolves"),
    @("Malik Monk", "PG,SG,SF", "Sacramento Kings"),
    @("Toumani Camara", "SF,PF", "Portland Trail Blazers"),
    @("Damian Lillard", "PG", "Milwaukee Bucks"),
    @("Isaiah Hartenstein", "C", "Oklahoma City Thunder"),
    @("Alex Caruso", "SG,SF", "Oklahoma City Thunder"),
    @("Brandon Miller", "SG,SF,PF", "Charlotte Hornets"),
    @("Brandon Ingram", "SG,SF,PF", "New Orleans Pelicans"),
    @("LaMelo Ball", "PG,SG", "Charlotte Hornets")
)

$r = 2
foreach ($row in $rows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $r = $r + 1
}
